# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '60.843.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.38%  '
$ws.Range('D3').Value = "'" + '2.454.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.92%  '
$ws.Range('D5').Value = "'" + '548.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.11%  '
$ws.Range('D6').Value = "'" + '145.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.33%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -3.75%  '
$ws.Range('D9').Value = "'" + '2.451.93'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.94%  '
$ws.Range('E10').Value = '  -8.21%  '
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = "'" + '5.36'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.01%  '
$ws.Range('D13').Value = "'" + '0.352'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.13%  '
$ws.Range('D14').Value = "'" + '26.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.79%  '
$ws.Range('D15').Value = "'" + '2.897.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.76%  '
$ws.Range('E16').Value = '  -9.02%  '
$ws.Range('D17').Value = "'" + '60.738.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('D18').Value = "'" + '2.455.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '11.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.47%  '
$ws.Range('E20').Value = '  -7.49%  '
$ws.Range('E21').Value = '  -7.08%  '
$ws.Range('D22').Value = "'" + '318.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.56%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'" + '63.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.53%  '
$ws.Range('E25').Value = '  -2.34%  '
$ws.Range('D26').Value = "'" + '0.0₃0975'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.37%  '
$ws.Range('D27').Value = "'" + '2.579.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.63%  '
$ws.Range('D28').Value = "'" + '539.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.77%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  -4.32%  '
$ws.Range('D31').Value = "'" + '8.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.22%  '
$ws.Range('D32').Value = "'" + '7.54'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.42%  '
$ws.Range('D33').Value = "'" + '0.149'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.83%  '
$ws.Range('D34').Value = "'" + '1.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.76%  '
$ws.Range('D35').Value = "'" + '1.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.61%  '
$ws.Range('D36').Value = "'" + '5.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -10.87%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = "'" + '4.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.37%  '
$ws.Range('E39').Value = '  -6.25%  '
$ws.Range('D40').Value = "'" + '18.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.95%  '
$ws.Range('D41').Value = "'" + '145.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E43').Value = '  -8.75%  '
$ws.Range('D44').Value = "'" + '39.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.15%  '
$ws.Range('E45').Value = '  -8.47%  '
$ws.Range('D46').Value = "'" + '146.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.62%  '
$ws.Range('E47').Value = '  -7.84%  '
$ws.Range('D48').Value = "'" + '20.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.58%  '
$ws.Range('E49').Value = '  -8.96%  '
$ws.Range('E50').Value = '  -7.38%  '
$ws.Range('D51').Value = "'" + '0.0937'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.08%  '
